$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape
$ws.Range("D2").Value = "90.618.44"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "3.101.64"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'231.47"
$ws.Range("E5").Value = "  +5.27%  "
$ws.Range("D6").Value = "'624.60"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'1.12"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("D8").Value = "'0.362"
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "3.105.90"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").Value = "'0.728"
$ws.Range("E11").Value = "  -5.40%  "
$ws.Range("D12").Value = "'0.197"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "'36.42"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "'5.48"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "90.599.87"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "3.671.47"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").Value = "3.120.73"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "'3.78"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'14.09"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("D21").Value = "'0.0000208"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").Value = "'440.56"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").Value = "'5.55"
$ws.Range("E23").Value = "  +6.35%  "
$ws.Range("D24").Value = "'8.88"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "'5.85"
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").Value = "'89.32"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").Value = "'12.22"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'9.45"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").Value = "'0.204"
$ws.Range("E32").Value = "  +17.34%  "
$ws.Range("D33").Value = "'26.40"
$ws.Range("E33").Value = "  +5.46%  "
$ws.Range("D34").Value = "'0.893"
$ws.Range("E34").Value = "  -10.43%  "
$ws.Range("D35").Value = "'0.150"
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("D36").Value = "'3.76"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "'508.58"
$ws.Range("E37").Value = "  -3.39%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").Value = "'1.91"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'7.03"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("D41").Value = "'0.0890"
$ws.Range("E41").Value = "  +4.59%  "
$ws.Range("D42").Value = "'0.410"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D45").Value = "'3.37"
$ws.Range("E45").Value = "  +49.99%  "
$ws.Range("D46").Value = "'1.90"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "'150.96"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").Value = "'0.685"
$ws.Range("E48").Value = "  +5.46%  "
$ws.Range("D49").Value = "'45.11"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("D50").Value = "'1.33"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").Value = "'4.43"
$ws.Range("E51").Value = "  +1.36%  "
